$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data row (row 2) with new case values and renumbered candidate id.
$ws.Range("A2").Value = "OBHCz204"
$ws.Range("B2").Value = 23090528
$ws.Range("C2").Value = "vlgxpgw72"
$ws.Range("D2").Value = "m6V!n5N#"
$ws.Range("F2").Value = "EigazUKb"
$ws.Range("G2").Value = "ReAG"
